$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = [double]"1.433944"
$ws.Range("H2").Value = [double]"4.301832"
$ws.Range("I2").Value = [double]"0.003882998715548277"
$ws.Range("J2").Value = [double]"0.003886188105009087"
$ws.Range("M2").Value = [double]"0.5804443333333333"
$ws.Range("N2").Value = [double]"1.741333"
$ws.Range("O2").Value = [double]"0.002431273010151717"
$ws.Range("P2").Value = [double]"0.002435427107574628"
$ws.Range("Q2").Value = [double]"0.8323246691173334"
$ws.Range("R2").Value = [double]"7.490922022056"
$ws.Range("S2").Value = [double]"9.440629975566313E-06"
$ws.Range("T2").Value = [double]"9.464527856073207E-06"
$ws.Range("G3").Value = [double]"1.433944"
$ws.Range("H3").Value = [double]"4.301832"
$ws.Range("I3").Value = [double]"0.003882998715548277"
$ws.Range("J3").Value = [double]"0.003886188105009087"
$ws.Range("O3").Value = [double]"0.0004752041289926495"
$ws.Range("P3").Value = [double]"0.00047601606752829"
$ws.Range("Q3").Value = [double]"0.1626819027626667"
$ws.Range("R3").Value = [double]"1.464137124864"
$ws.Range("S3").Value = [double]"1.845217022501696E-06"
$ws.Range("T3").Value = [double]"1.849887979421643E-06"
$ws.Range("G4").Value = [double]"1.433944"
$ws.Range("H4").Value = [double]"4.301832"
$ws.Range("I4").Value = [double]"0.003882998715548277"
$ws.Range("J4").Value = [double]"0.003886188105009087"
$ws.Range("M4").Value = [double]"136.1000366666667"
$ws.Range("N4").Value = [double]"408.30011"
$ws.Range("O4").Value = [double]"0.5700742118164518"
$ws.Range("P4").Value = [double]"0.5710482463260632"
$ws.Range("Q4").Value = [double]"195.1598309779467"
$ws.Range("R4").Value = [double]"1756.43847880152"
$ws.Range("S4").Value = [double]"0.002213597432250479"
$ws.Range("T4").Value = [double]"0.002219200902258646"
$ws.Range("G5").Value = [double]"1.433944"
$ws.Range("H5").Value = [double]"4.301832"
$ws.Range("I5").Value = [double]"0.003882998715548277"
$ws.Range("J5").Value = [double]"0.003886188105009087"
$ws.Range("M5").Value = [double]"1.221658"
$ws.Range("N5").Value = [double]"2.443316"
$ws.Range("O5").Value = [double]"0.005117086949542552"
$ws.Range("P5").Value = [double]"0.003417220037046797"
$ws.Range("Q5").Value = [double]"1.751789159152"
$ws.Range("R5").Value = [double]"10.510734954912"
$ws.Range("S5").Value = [double]"1.986964205242258E-05"
$ws.Range("T5").Value = [double]"1.327995986016997E-05"
$ws.Range("G6").Value = [double]"1.433944"
$ws.Range("H6").Value = [double]"4.301832"
$ws.Range("I6").Value = [double]"0.003882998715548277"
$ws.Range("J6").Value = [double]"0.003886188105009087"
$ws.Range("M6").Value = [double]"100.7253213333333"
$ws.Range("N6").Value = [double]"302.175964"
$ws.Range("O6").Value = [double]"0.4219022240948613"
$ws.Range("P6").Value = [double]"0.4226230904617871"
$ws.Range("Q6").Value = [double]"144.4344701740054"
$ws.Range("R6").Value = [double]"1299.910231566048"
$ws.Range("S6").Value = [double]"0.001638245794247308"
$ws.Range("T6").Value = [double]"0.001642392827054776"
$ws.Range("I7").Value = [double]"0.7877262822264709"
$ws.Range("J7").Value = [double]"0.7883732991550308"
$ws.Range("M7").Value = [double]"0.5804443333333333"
$ws.Range("N7").Value = [double]"1.741333"
$ws.Range("O7").Value = [double]"0.002431273010151717"
$ws.Range("P7").Value = [double]"0.002435427107574628"
$ws.Range("Q7").Value = [double]"168.8499186424784"
$ws.Range("R7").Value = [double]"1519.649267782305"
$ws.Range("S7").Value = [double]"0.001915177649364373"
$ws.Range("T7").Value = [double]"0.001920025703650204"
$ws.Range("I8").Value = [double]"0.7877262822264709"
$ws.Range("J8").Value = [double]"0.7883732991550308"
$ws.Range("O8").Value = [double]"0.0004752041289926495"
$ws.Range("P8").Value = [double]"0.00047601606752829"
$ws.Range("S8").Value = [double]"0.0003743307818300481"
$ws.Range("T8").Value = [double]"0.000375278357608082"
$ws.Range("I9").Value = [double]"0.7877262822264709"
$ws.Range("J9").Value = [double]"0.7883732991550308"
$ws.Range("M9").Value = [double]"136.1000366666667"
$ws.Range("N9").Value = [double]"408.30011"
$ws.Range("O9").Value = [double]"0.5700742118164518"
$ws.Range("P9").Value = [double]"0.5710482463260632"
$ws.Range("Q9").Value = [double]"39591.18695574883"
$ws.Range("R9").Value = [double]"356320.6826017394"
$ws.Range("S9").Value = [double]"0.4490624394673592"
$ws.Range("T9").Value = [double]"0.4501991899327732"
$ws.Range("I10").Value = [double]"0.7877262822264709"
$ws.Range("J10").Value = [double]"0.7883732991550308"
$ws.Range("M10").Value = [double]"1.221658"
$ws.Range("N10").Value = [double]"2.443316"
$ws.Range("O10").Value = [double]"0.005117086949542552"
$ws.Range("P10").Value = [double]"0.003417220037046797"
$ws.Range("Q10").Value = [double]"355.37749627831"
$ws.Range("R10").Value = [double]"2132.26497766986"
$ws.Range("S10").Value = [double]"0.004030863878592747"
$ws.Range("T10").Value = [double]"0.00269404503454526"
$ws.Range("I11").Value = [double]"0.7877262822264709"
$ws.Range("J11").Value = [double]"0.7883732991550308"
$ws.Range("M11").Value = [double]"100.7253213333333"
$ws.Range("N11").Value = [double]"302.175964"
$ws.Range("O11").Value = [double]"0.4219022240948613"
$ws.Range("P11").Value = [double]"0.4226230904617871"
$ws.Range("Q11").Value = [double]"29300.763804001"
$ws.Range("R11").Value = [double]"263706.874236009"
$ws.Range("S11").Value = [double]"0.3323434704493244"
$ws.Range("T11").Value = [double]"0.3331847601264541"
$ws.Range("G12").Value = [double]"54.70735966666666"
$ws.Range("H12").Value = [double]"164.122079"
$ws.Range("I12").Value = [double]"0.1481428893434501"
$ws.Range("J12").Value = [double]"0.1482645698807303"
$ws.Range("M12").Value = [double]"0.5804443333333333"
$ws.Range("N12").Value = [double]"1.741333"
$ws.Range("O12").Value = [double]"0.002431273010151717"
$ws.Range("P12").Value = [double]"0.002435427107574628"
$ws.Range("Q12").Value = [double]"31.75457691014522"
$ws.Range("R12").Value = [double]"285.791192191307"
$ws.Range("S12").Value = [double]"0.0003601758085066228"
$ws.Range("T12").Value = [double]"0.0003610875525804233"
$ws.Range("G13").Value = [double]"54.70735966666666"
$ws.Range("H13").Value = [double]"164.122079"
$ws.Range("I13").Value = [double]"0.1481428893434501"
$ws.Range("J13").Value = [double]"0.1482645698807303"
$ws.Range("O13").Value = [double]"0.0004752041289926495"
$ws.Range("P13").Value = [double]"0.00047601606752829"
$ws.Range("Q13").Value = [double]"6.206586425756443"
$ws.Range("R13").Value = [double]"55.85927783180799"
$ws.Range("S13").Value = [double]"7.039811269690868E-05"
$ws.Range("T13").Value = [double]"7.057631750839857E-05"
$ws.Range("G14").Value = [double]"54.70735966666666"
$ws.Range("H14").Value = [double]"164.122079"
$ws.Range("I14").Value = [double]"0.1481428893434501"
$ws.Range("J14").Value = [double]"0.1482645698807303"
$ws.Range("M14").Value = [double]"136.1000366666667"
$ws.Range("N14").Value = [double]"408.30011"
$ws.Range("O14").Value = [double]"0.5700742118164518"
$ws.Range("P14").Value = [double]"0.5710482463260632"
$ws.Range("Q14").Value = [double]"7445.673656569855"
$ws.Range("R14").Value = [double]"67011.06290912868"
$ws.Range("S14").Value = [double]"0.08445244087867917"
$ws.Range("T14").Value = [double]"0.08466622262267906"
$ws.Range("G15").Value = [double]"54.70735966666666"
$ws.Range("H15").Value = [double]"164.122079"
$ws.Range("I15").Value = [double]"0.1481428893434501"
$ws.Range("J15").Value = [double]"0.1482645698807303"
$ws.Range("M15").Value = [double]"1.221658"
$ws.Range("N15").Value = [double]"2.443316"
$ws.Range("O15").Value = [double]"0.005117086949542552"
$ws.Range("P15").Value = [double]"0.003417220037046797"
$ws.Range("Q15").Value = [double]"66.83368359566066"
$ws.Range("R15").Value = [double]"401.002101573964"
$ws.Range("S15").Value = [double]"0.0007580600457268951"
$ws.Range("T15").Value = [double]"0.0005066526589805564"
$ws.Range("G16").Value = [double]"54.70735966666666"
$ws.Range("H16").Value = [double]"164.122079"
$ws.Range("I16").Value = [double]"0.1481428893434501"
$ws.Range("J16").Value = [double]"0.1482645698807303"
$ws.Range("M16").Value = [double]"100.7253213333333"
$ws.Range("N16").Value = [double]"302.175964"
$ws.Range("O16").Value = [double]"0.4219022240948613"
$ws.Range("P16").Value = [double]"0.4226230904617871"
$ws.Range("Q16").Value = [double]"5510.416381723239"
$ws.Range("R16").Value = [double]"49593.74743550915"
$ws.Range("S16").Value = [double]"0.06250181449784054"
$ws.Range("T16").Value = [double]"0.06266003072898181"
$ws.Range("G17").Value = [double]"0.909222"
$ws.Range("H17").Value = [double]"1.818444"
$ws.Range("I17").Value = [double]"0.002462096049879378"
$ws.Range("J17").Value = [double]"0.001642745565709015"
$ws.Range("M17").Value = [double]"0.5804443333333333"
$ws.Range("N17").Value = [double]"1.741333"
$ws.Range("O17").Value = [double]"0.002431273010151717"
$ws.Range("P17").Value = [double]"0.002435427107574628"
$ws.Range("Q17").Value = [double]"0.527752757642"
$ws.Range("R17").Value = [double]"3.166516545852"
$ws.Range("S17").Value = [double]"5.986027674472889E-06"
$ws.Range("T17").Value = [double]"4.000787081575754E-06"
$ws.Range("G18").Value = [double]"0.909222"
$ws.Range("H18").Value = [double]"1.818444"
$ws.Range("I18").Value = [double]"0.002462096049879378"
$ws.Range("J18").Value = [double]"0.001642745565709015"
$ws.Range("O18").Value = [double]"0.0004752041289926495"
$ws.Range("P18").Value = [double]"0.00047601606752829"
$ws.Range("Q18").Value = [double]"0.103151842048"
$ws.Range("R18").Value = [double]"0.6189110522879999"
$ws.Range("S18").Value = [double]"1.169998208879173E-06"
$ws.Range("T18").Value = [double]"7.819732841383416E-07"
$ws.Range("G19").Value = [double]"0.909222"
$ws.Range("H19").Value = [double]"1.818444"
$ws.Range("I19").Value = [double]"0.002462096049879378"
$ws.Range("J19").Value = [double]"0.001642745565709015"
$ws.Range("M19").Value = [double]"136.1000366666667"
$ws.Range("N19").Value = [double]"408.30011"
$ws.Range("O19").Value = [double]"0.5700742118164518"
$ws.Range("P19").Value = [double]"0.5710482463260632"
$ws.Range("Q19").Value = [double]"123.74514753814"
$ws.Range("R19").Value = [double]"742.47088522884"
$ws.Range("S19").Value = [double]"0.001403577465051386"
$ws.Range("T19").Value = [double]"0.0009380869744580497"
$ws.Range("G20").Value = [double]"0.909222"
$ws.Range("H20").Value = [double]"1.818444"
$ws.Range("I20").Value = [double]"0.002462096049879378"
$ws.Range("J20").Value = [double]"0.001642745565709015"
$ws.Range("M20").Value = [double]"1.221658"
$ws.Range("N20").Value = [double]"2.443316"
$ws.Range("O20").Value = [double]"0.005117086949542552"
$ws.Range("P20").Value = [double]"0.003417220037046797"
$ws.Range("Q20").Value = [double]"1.110758330076"
$ws.Range("R20").Value = [double]"4.443033320303999"
$ws.Range("S20").Value = [double]"1.259875956535804E-05"
$ws.Range("T20").Value = [double]"5.613623062910622E-06"
$ws.Range("G21").Value = [double]"0.909222"
$ws.Range("H21").Value = [double]"1.818444"
$ws.Range("I21").Value = [double]"0.002462096049879378"
$ws.Range("J21").Value = [double]"0.001642745565709015"
$ws.Range("M21").Value = [double]"100.7253213333333"
$ws.Range("N21").Value = [double]"302.175964"
$ws.Range("O21").Value = [double]"0.4219022240948613"
$ws.Range("P21").Value = [double]"0.4226230904617871"
$ws.Range("Q21").Value = [double]"91.58167811333601"
$ws.Range("R21").Value = [double]"549.4900686800161"
$ws.Range("S21").Value = [double]"0.001038763799379282"
$ws.Range("T21").Value = [double]"0.0006942622078223407"
$ws.Range("G22").Value = [double]"21.33956566666667"
$ws.Range("H22").Value = [double]"64.018697"
$ws.Range("I22").Value = [double]"0.05778573366465133"
$ws.Range("J22").Value = [double]"0.05783319729352075"
$ws.Range("M22").Value = [double]"0.5804443333333333"
$ws.Range("N22").Value = [double]"1.741333"
$ws.Range("O22").Value = [double]"0.002431273010151717"
$ws.Range("P22").Value = [double]"0.002435427107574628"
$ws.Range("Q22").Value = [double]"12.38642996701122"
$ws.Range("R22").Value = [double]"111.477869703101"
$ws.Range("S22").Value = [double]"0.0001404928946306823"
$ws.Range("T22").Value = [double]"0.0001408485364063521"
$ws.Range("G23").Value = [double]"21.33956566666667"
$ws.Range("H23").Value = [double]"64.018697"
$ws.Range("I23").Value = [double]"0.05778573366465133"
$ws.Range("J23").Value = [double]"0.05783319729352075"
$ws.Range("O23").Value = [double]"0.0004752041289926495"
$ws.Range("P23").Value = [double]"0.00047601606752829"
$ws.Range("Q23").Value = [double]"2.420987951260444"
$ws.Range("R23").Value = [double]"21.788891561344"
$ws.Range("S23").Value = [double]"2.746001923431185E-05"
$ws.Range("T23").Value = [double]"2.752953114824949E-05"
$ws.Range("G24").Value = [double]"21.33956566666667"
$ws.Range("H24").Value = [double]"64.018697"
$ws.Range("I24").Value = [double]"0.05778573366465133"
$ws.Range("J24").Value = [double]"0.05783319729352075"
$ws.Range("M24").Value = [double]"136.1000366666667"
$ws.Range("N24").Value = [double]"408.30011"
$ws.Range("O24").Value = [double]"0.5700742118164518"
$ws.Range("P24").Value = [double]"0.5710482463260632"
$ws.Range("Q24").Value = [double]"2904.315669684075"
$ws.Range("R24").Value = [double]"26138.84102715667"
$ws.Range("S24").Value = [double]"0.03294215657311151"
$ws.Range("T24").Value = [double]"0.03302554589389425"
$ws.Range("G25").Value = [double]"21.33956566666667"
$ws.Range("H25").Value = [double]"64.018697"
$ws.Range("I25").Value = [double]"0.05778573366465133"
$ws.Range("J25").Value = [double]"0.05783319729352075"
$ws.Range("M25").Value = [double]"1.221658"
$ws.Range("N25").Value = [double]"2.443316"
$ws.Range("O25").Value = [double]"0.005117086949542552"
$ws.Range("P25").Value = [double]"0.003417220037046797"
$ws.Range("Q25").Value = [double]"26.06965111320867"
$ws.Range("R25").Value = [double]"156.417906679252"
$ws.Range("S25").Value = [double]"0.000295694623605129"
$ws.Range("T25").Value = [double]"0.0001976287605978997"
$ws.Range("G26").Value = [double]"21.33956566666667"
$ws.Range("H26").Value = [double]"64.018697"
$ws.Range("I26").Value = [double]"0.05778573366465133"
$ws.Range("J26").Value = [double]"0.05783319729352075"
$ws.Range("M26").Value = [double]"100.7253213333333"
$ws.Range("N26").Value = [double]"302.175964"
$ws.Range("O26").Value = [double]"0.4219022240948613"
$ws.Range("P26").Value = [double]"0.4226230904617871"
$ws.Range("Q26").Value = [double]"2149.434608888768"
$ws.Range("R26").Value = [double]"19344.91147999891"
$ws.Range("S26").Value = [double]"0.02437992955406969"
$ws.Range("T26").Value = [double]"0.024441644571474"
